# Daily attendance processing - normalize "Recorded By" (column G) ordering.
#
# Rule (derived from the target diff): within each comma-separated list of
# recorders, stably reorder so that:
#   - the exact (case-sensitive) token "System" always sorts last, and
#   - the token "admin@admin.com" sorts after every other (non-"System")
#     recorder, but still before "System";
#   - all other tokens (e.g. the lowercase "system", "backup@backdoor.com",
#     "dnasr281@gmail.com") keep their original relative order.
#
# NOTE: this engine's -eq/-ceq/-cmatch string operators are all
# case-INsensitive, so an explicit char-code comparison is used to tell
# "System" apart from "system".

function Test-ExactEquals($a, $b) {
    if ($a.Length -ne $b.Length) {
        return $false
    }
    for ($i = 0; $i -lt $a.Length; $i++) {
        $ca = [int][char]$a.Substring($i, 1)
        $cb = [int][char]$b.Substring($i, 1)
        if ($ca -ne $cb) {
            return $false
        }
    }
    return $true
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $v = $cell.Value2

    if ($v -eq $null) {
        continue
    }
    if ($v -eq "") {
        continue
    }

    $parts = $v -split ", "

    $primary = @()
    $admin = @()
    $system = @()

    foreach ($p in $parts) {
        if (Test-ExactEquals $p "System") {
            $system += $p
        } elseif ($p -eq "admin@admin.com") {
            $admin += $p
        } else {
            $primary += $p
        }
    }

    $result = $primary + $admin + $system
    $joined = $result -join ", "

    if ($joined -ne $v) {
        $cell.Value2 = $joined
    }
}
